$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest crypto data snapshot.
# D-column price values are forced to Text format first, since many of them
# look like plain numbers (e.g. "605.24") and Excel would otherwise silently
# convert them to floating point numbers, losing the original text formatting
# (trailing zeros, thousand-dot grouping, etc.). The style is reset back to
# "Normal" afterwards so no stray cell formatting is introduced.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.409.85'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.52%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.647.85'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.69%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '605.24'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '155.99'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.98%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  -0.59%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.645.36'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.62%  '
$ws.Range('E10').Value = '  +7.45%  '
$ws.Range('E11').Value = '  +0.96%  '
$ws.Range('E13').Value = '  +1.57%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '29.90'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.30%  '
$ws.Range('E15').Value = '  +14.07%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.126.15'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.85%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.184.67'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.36%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.647.77'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.71'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.33%  '
$ws.Range('E20').Value = '  +2.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '358.21'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.41%  '
$ws.Range('E22').Value = '  +4.99%  '
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.64'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.84%  '
$ws.Range('E25').Value = '  -0.89%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.41'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.32%  '
$ws.Range('E27').Value = '  +14.50%  '
$ws.Range('E29').Value = '  +2.27%  '
$ws.Range('E30').Value = '  -2.58%  '
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.16'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.87%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '521.35'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.77'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.89%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.50'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.33'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E37').Value = '  +2.19%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '20.67'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '162.47'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.79%  '
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('E41').Value = '  -1.85%  '
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('E43').Value = '  +0.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '165.42'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.09%  '
$ws.Range('E45').Value = '  -0.20%  '
$ws.Range('E46').Value = '  +3.67%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '22.89'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.85%  '
$ws.Range('E49').Value = '  +3.86%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.649'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.33%  '
$ws.Range('E51').Value = '  +0.08%  '
